$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: 'ALC'
$ws.Range("H15").Value = 975.25
$ws.Range("I15").Value = 975.25
$ws.Range("K15").Value = 2925.75
$ws.Range("M15").Value = -2756.75

# Row 96: 'ALC'
$ws.Range("H96").Value = 2728.75
$ws.Range("I96").Value = 5026
$ws.Range("K96").Value = 15078
$ws.Range("M96").Value = -13705

# Row 98: 'ALC'
$ws.Range("H98").Value = 4142.5356
$ws.Range("I98").Value = 4152.737
$ws.Range("J98").Value = 4121
$ws.Range("K98").Value = 4152.737
$ws.Range("L98").Value = 4121
$ws.Range("M98").Value = -2654.737
$ws.Range("N98").Value = -7117

# Row 122: 'ALC'
$ws.Range("H122").Value = 4142.5356
$ws.Range("I122").Value = 4152.737
$ws.Range("J122").Value = 4121
$ws.Range("K122").Value = 12458.211
$ws.Range("L122").Value = 12363
$ws.Range("M122").Value = -10008.211
$ws.Range("N122").Value = -17263

# Row 137: 'ALC'
$ws.Range("H137").Value = 34348.867
$ws.Range("I137").Value = 672.5263
$ws.Range("K137").Value = 2017.5789
$ws.Range("M137").Value = 532.4211

# Row 138: 'ALC'
$ws.Range("H138").Value = 2395.739
$ws.Range("I138").Value = 2620.4324
$ws.Range("J138").Value = 2244.5818
$ws.Range("K138").Value = 7861.297200000001
$ws.Range("L138").Value = 6733.7454
$ws.Range("M138").Value = -2721.297200000001
$ws.Range("N138").Value = -17013.7454

$ws = $wb.Worksheets.Item("ARM")
# Row 32: 'ARM'
$ws.Range("H32").Value = 5518.6885
$ws.Range("I32").Value = 3048.635
$ws.Range("K32").Value = 3048.635
$ws.Range("M32").Value = -2761.635

# Row 74: 'ARM'
$ws.Range("H74").Value = 744.075
$ws.Range("I74").Value = 547.2
$ws.Range("K74").Value = 547.2
$ws.Range("M74").Value = 326.8

# Row 77: 'ARM'
$ws.Range("H77").Value = 744.075
$ws.Range("I77").Value = 547.2
$ws.Range("K77").Value = 2736
$ws.Range("M77").Value = 1632

# Row 102: 'ARM'
$ws.Range("H102").Value = 1877.5
$ws.Range("I102").Value = 1877.5
$ws.Range("K102").Value = 1877.5
$ws.Range("M102").Value = -255.5

# Row 122: 'ARM'
$ws.Range("H122").Value = 1144
$ws.Range("I122").Value = 1154.3334
$ws.Range("K122").Value = 3463.0002
$ws.Range("M122").Value = -1013.0002

# Row 132: 'ARM'
$ws.Range("H132").Value = 2104.1936
$ws.Range("I132").Value = 1651.8823
$ws.Range("J132").Value = 2653.4285
$ws.Range("K132").Value = 4955.6469
$ws.Range("L132").Value = 7960.2855
$ws.Range("M132").Value = -2425.6469
$ws.Range("N132").Value = -13020.2855

# Row 141: 'ARM'
$ws.Range("H141").Value = 27000
$ws.Range("J141").Value = 27000
$ws.Range("L141").Value = 27000
$ws.Range("N141").Value = -37360

$ws = $wb.Worksheets.Item("BSM")
# Row 86: 'BSM'
$ws.Range("H86").Value = 2035.6666
$ws.Range("I86").Value = 1554.25
$ws.Range("J86").Value = 2998.5
$ws.Range("K86").Value = 1554.25
$ws.Range("L86").Value = 2998.5
$ws.Range("M86").Value = -431.25
$ws.Range("N86").Value = -5244.5

# Row 89: 'BSM'
$ws.Range("H89").Value = 2035.6666
$ws.Range("I89").Value = 1554.25
$ws.Range("J89").Value = 2998.5
$ws.Range("K89").Value = 7771.25
$ws.Range("L89").Value = 14992.5
$ws.Range("M89").Value = -2155.25
$ws.Range("N89").Value = -26224.5

# Row 94: 'BSM'
$ws.Range("H94").Value = 797.625
$ws.Range("I94").Value = 401.75
$ws.Range("J94").Value = 2777
$ws.Range("K94").Value = 401.75
$ws.Range("L94").Value = 2777
$ws.Range("M94").Value = 49.25
$ws.Range("N94").Value = -3679

# Row 137: 'BSM'
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

# Row 140: 'BSM'
$ws.Range("H140").Value = 46099.875
$ws.Range("J140").Value = 46099.875
$ws.Range("L140").Value = 46099.875
$ws.Range("N140").Value = -56459.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 'CRP'
$ws.Range("H31").Value = 1821.1316
$ws.Range("I31").Value = 1487.2693
$ws.Range("K31").Value = 1487.2693
$ws.Range("M31").Value = -1192.2693

# Row 34: 'CRP'
$ws.Range("H34").Value = 1821.1316
$ws.Range("I34").Value = 1487.2693
$ws.Range("K34").Value = 1487.2693
$ws.Range("M34").Value = -1285.2693

# Row 58: 'CRP'
$ws.Range("H58").Value = 3107722
$ws.Range("I58").Value = 6213729
$ws.Range("K58").Value = 6213729
$ws.Range("M58").Value = -6213526

# Row 99: 'CRP'
$ws.Range("H99").Value = 2635.5715
$ws.Range("J99").Value = 2595.6
$ws.Range("L99").Value = 2595.6
$ws.Range("N99").Value = -5591.6

# Row 126: 'CRP'
$ws.Range("H126").Value = 2635.5715
$ws.Range("J126").Value = 2595.6
$ws.Range("L126").Value = 7786.799999999999
$ws.Range("N126").Value = -12726.8

# Row 136: 'CRP'
$ws.Range("H136").Value = 3107722
$ws.Range("I136").Value = 6213729
$ws.Range("K136").Value = 18641187
$ws.Range("M136").Value = -18638637

$ws = $wb.Worksheets.Item("CUL")
# Row 80: 'CUL'
$ws.Range("H80").Value = 1915.8334
$ws.Range("J80").Value = 2250
$ws.Range("L80").Value = 6750
$ws.Range("N80").Value = -8622

# Row 83: 'CUL'
$ws.Range("H83").Value = 1915.8334
$ws.Range("J83").Value = 2250
$ws.Range("L83").Value = 20250
$ws.Range("N83").Value = -29610

# Row 119: 'CUL'
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

# Row 131: 'CUL'
$ws.Range("H131").Value = 15260.903
$ws.Range("I131").Value = 269.66666
$ws.Range("J131").Value = 16867.107
$ws.Range("K131").Value = 808.9999799999999
$ws.Range("L131").Value = 50601.321
$ws.Range("M131").Value = 4231.00002
$ws.Range("N131").Value = -60681.321

# Row 136: 'CUL'
$ws.Range("H136").Value = 1430.7894
$ws.Range("I136").Value = 1361.5625
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 4084.6875
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = 1015.3125
$ws.Range("N136").Value = -15600

# Row 140: 'CUL'
$ws.Range("H140").Value = 2404.9048
$ws.Range("I140").Value = 1733.6111
$ws.Range("J140").Value = 6432.6665
$ws.Range("K140").Value = 5200.8333
$ws.Range("L140").Value = 19297.9995
$ws.Range("M140").Value = -20.83330000000024
$ws.Range("N140").Value = -29657.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 122: 'GSM'
$ws.Range("H122").Value = 1330
$ws.Range("I122").Value = 1228.6364
$ws.Range("K122").Value = 3685.9092
$ws.Range("M122").Value = -1235.9092

# Row 126: 'GSM'
$ws.Range("H126").Value = 2359346.2
$ws.Range("I126").Value = 2528218.5
$ws.Range("K126").Value = 7584655.5
$ws.Range("M126").Value = -7582185.5

# Row 132: 'GSM'
$ws.Range("H132").Value = 1834918.4
$ws.Range("I132").Value = 2406803.8
$ws.Range("J132").Value = 4885.2
$ws.Range("K132").Value = 7220411.399999999
$ws.Range("L132").Value = 14655.6
$ws.Range("M132").Value = -7217881.399999999
$ws.Range("N132").Value = -19715.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7: 'LTW'
$ws.Range("H7").Value = 3183.4
$ws.Range("I7").Value = 4238.3335
$ws.Range("J7").Value = 2731.2856
$ws.Range("K7").Value = 4238.3335
$ws.Range("L7").Value = 2731.2856
$ws.Range("M7").Value = -4126.3335
$ws.Range("N7").Value = -2955.2856

# Row 32: 'LTW'
$ws.Range("H32").Value = 8180
$ws.Range("I32").Value = 4800
$ws.Range("K32").Value = 4800
$ws.Range("M32").Value = -4483

# Row 40: 'LTW'
$ws.Range("H40").Value = 9445.857
$ws.Range("J40").Value = 8158.8
$ws.Range("L40").Value = 8158.8
$ws.Range("N40").Value = -8430.799999999999

# Row 48: 'LTW'
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0

# Row 55: 'LTW'
$ws.Range("H55").Value = 583.3889
$ws.Range("I55").Value = 498
$ws.Range("K55").Value = 498
$ws.Range("M55").Value = -325

# Row 61: 'LTW'
$ws.Range("H61").Value = 2964.2222
$ws.Range("I61").Value = 2613.1667
$ws.Range("J61").Value = 3666.3333
$ws.Range("K61").Value = 2613.1667
$ws.Range("L61").Value = 3666.3333
$ws.Range("M61").Value = -2411.1667
$ws.Range("N61").Value = -4070.3333

# Row 93: 'LTW'
$ws.Range("H93").Value = 1119.4286
$ws.Range("J93").Value = 3122.75
$ws.Range("L93").Value = 3122.75
$ws.Range("N93").Value = -5618.75

# Row 113: 'LTW'
$ws.Range("H113").Value = 2964.2222
$ws.Range("I113").Value = 2613.1667
$ws.Range("J113").Value = 3666.3333
$ws.Range("K113").Value = 2613.1667
$ws.Range("L113").Value = 3666.3333
$ws.Range("M113").Value = -443.1667000000002
$ws.Range("N113").Value = -8006.3333

# Row 126: 'LTW'
$ws.Range("H126").Value = 3183.4
$ws.Range("I126").Value = 4238.3335
$ws.Range("J126").Value = 2731.2856
$ws.Range("K126").Value = 12715.0005
$ws.Range("L126").Value = 8193.856800000001
$ws.Range("M126").Value = -10245.0005
$ws.Range("N126").Value = -13133.8568

# Row 132: 'LTW'
$ws.Range("H132").Value = 2330.0293
$ws.Range("I132").Value = 1802.25
$ws.Range("J132").Value = 2617.9092
$ws.Range("K132").Value = 5406.75
$ws.Range("L132").Value = 7853.7276
$ws.Range("M132").Value = -2876.75
$ws.Range("N132").Value = -12913.7276

# Row 133: 'LTW'
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -94060

$ws = $wb.Worksheets.Item("WVR")
# Row 96: 'WVR'
$ws.Range("H96").Value = 1677.2778
$ws.Range("I96").Value = 1542.6666
$ws.Range("J96").Value = 1811.8889
$ws.Range("K96").Value = 1542.6666
$ws.Range("L96").Value = 1811.8889
$ws.Range("M96").Value = -169.6666
$ws.Range("N96").Value = -4557.8889

# Row 132: 'WVR'
$ws.Range("H132").Value = 1219.88
$ws.Range("I132").Value = 1070.6842
$ws.Range("K132").Value = 3212.0526
$ws.Range("M132").Value = -682.0526

# Row 133: 'WVR'
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
